$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Unit Processes")
$ws2 = $wb.Worksheets.Item("Fuels")

# --- Fuels sheet header row updates -------------------------------------
# D1: "biogenic" -> "H2O"; add new E1 header "meta-notes"
$ws2.Range("D1").Value = "H2O"
$ws2.Range("E1").Value = "meta-notes"

# --- New fuel row: steam ---------------------------------------------
$ws2.Range("A8").Value = "steam"
$ws2.Range("B8").Value = 2.77
$ws2.Range("C8").Value = 0
$ws2.Range("D8").Value = 1

# --- New fuel row: coke -------------------------------------------------
$ws2.Range("A9").Value = "coke"
$ws2.Range("B9").Value = 29.01
$ws2.Range("C9").Value = 3.23
$ws2.Range("D9").Value = 0

# Row 2 (units row): B2 "(mj/t)" -> "(gj/t)"; D2 "(t/t)" -> "(t/t combusted)"
$ws2.Range("B2").Value = "(gj/t)"
$ws2.Range("D2").Value = "(t/t combusted)"

# --- New fuel row: Eurofer electricity mix proxy -------------------------
$ws2.Range("A10").Value = "Eurofer electricity mix proxy"
$ws2.Range("B10").Value = 1
$ws2.Range("C10").Value = 0.11
$ws2.Range("D10").Value = 0

# --- Updated fuel data ----------------------------------------------------
# coal: CO2 factor 3.2 -> 3.19
$ws2.Range("C3").Value = 3.19

# charcoal: CO2 factor 1.8 -> 2.93 ; biogenic/H2O flag 1 -> 0
$ws2.Range("C4").Value = 2.93
$ws2.Range("D4").Value = 0

# --- Column A on Fuels sheet widened to fit new longer labels ----------
$ws2.Range("A1").ColumnWidth = 28.1666666666667

# --- Selection / active sheet: author was last working on Fuels sheet --
$ws1.Range("E4").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("G14").Select() | Out-Null
